# Apply the "Last version 7-12-2025 ready" edit:
#  - Shorten several answer strings on the History1 and Geography sheets.
#  - Move the active tab/selection back to Biology (first sheet), updating
#    each sheet's selection/active-cell state accordingly.

$wb = $excel.ActiveWorkbook

# --- History1 sheet: shorten two answers ---
$wsHistory = $wb.Worksheets.Item("History1")
$wsHistory.Range("B3").Value = "Berlin Wall"
$wsHistory.Range("B2").Value = "The Black Death"

# --- Geography sheet: shorten three answers ---
$wsGeography = $wb.Worksheets.Item("Geography")
$wsGeography.Range("B2").Value = "Nile"
$wsGeography.Range("B4").Value = "Istanbul"
$wsGeography.Range("B5").Value = "Pacific"

# --- Update selections on each sheet ---
$wsBiology = $wb.Worksheets.Item("Biology")

# History1 was previously the selected tab with a full-column selection;
# now it just has a regular cell selection and is no longer the active tab.
$wsHistory.Activate() | Out-Null
$wsHistory.Range("B10").Select() | Out-Null

# Biology becomes the active/selected tab with A12 selected.
$wsBiology.Activate() | Out-Null
$wsBiology.Range("A12").Select() | Out-Null
